# The edit renames the embedded logo picture shapes (wp:docPr / pic:cNvPr
# "name" attributes only -- their "descr" stays untouched) inside the two
# footers (Pearson logo: image1.png -> image2.png) and the two headers
# (BTec logo: image2.jpg -> image1.jpg).
#
# InlineShape has no writable "Name" property on the Word object model, so
# the rename is performed the way Word itself exposes raw-OOXML edits: pull
# the whole package as WordOpenXML, patch the two "name=" attribute values
# with a scoped Find/Replace, and write the package back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
